$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 2
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = 3
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = 1
